# Apply updated crypto price/volume figures to the worksheet.
# Generated from the authoritative cell-level diff of the commit
# "Updated cryptos list on Mon Jul 31 13:42:56 UTC 2023 with GitHub Actions".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{ Cell = 'D2'; Value = '29.435.75' },
    @{ Cell = 'E2'; Value = '  +0.23%  ' },
    @{ Cell = 'D3'; Value = '1.870.04' },
    @{ Cell = 'E3'; Value = '  -0.48%  ' },
    @{ Cell = 'E4'; Value = '  -0.12%  ' },
    @{ Cell = 'D5'; Value = '243.79' },
    @{ Cell = 'E5'; Value = '  +0.26%  ' },
    @{ Cell = 'D6'; Value = '0.7049' },
    @{ Cell = 'E6'; Value = '  -2.34%  ' },
    @{ Cell = 'E7'; Value = '  -0.10%  ' },
    @{ Cell = 'D8'; Value = '0.07924' },
    @{ Cell = 'E8'; Value = '  -1.32%  ' },
    @{ Cell = 'D9'; Value = '0.3134' },
    @{ Cell = 'E9'; Value = '  -0.12%  ' },
    @{ Cell = 'D10'; Value = '24.51' },
    @{ Cell = 'E10'; Value = '  -1.61%  ' },
    @{ Cell = 'D11'; Value = '0.07780' },
    @{ Cell = 'E11'; Value = '  -4.63%  ' },
    @{ Cell = 'D12'; Value = '1.860.88' },
    @{ Cell = 'E12'; Value = '  -1.06%  ' },
    @{ Cell = 'D13'; Value = '93.36' },
    @{ Cell = 'E13'; Value = '  -1.30%  ' },
    @{ Cell = 'D14'; Value = '5.158' },
    @{ Cell = 'E14'; Value = '  -1.32%  ' },
    @{ Cell = 'D15'; Value = '0.7002' },
    @{ Cell = 'E15'; Value = '  -1.50%  ' },
    @{ Cell = 'D16'; Value = '6.494' },
    @{ Cell = 'E16'; Value = '  +1.35%  ' },
    @{ Cell = 'D17'; Value = '0.000008543' },
    @{ Cell = 'E17'; Value = '  +0.77%  ' },
    @{ Cell = 'D18'; Value = '29.396.89' },
    @{ Cell = 'E18'; Value = '  +0.07%  ' },
    @{ Cell = 'D19'; Value = '251.27' },
    @{ Cell = 'E19'; Value = '  +2.86%  ' },
    @{ Cell = 'D20'; Value = '2.128.93' },
    @{ Cell = 'E20'; Value = '  +0.14%  ' },
    @{ Cell = 'D21'; Value = '13.07' },
    @{ Cell = 'E21'; Value = '  -1.58%  ' },
    @{ Cell = 'D22'; Value = '0.9997' },
    @{ Cell = 'E22'; Value = '  -0.21%  ' },
    @{ Cell = 'D23'; Value = '7.600' },
    @{ Cell = 'E23'; Value = '  -1.69%  ' },
    @{ Cell = 'E24'; Value = '  -0.02%  ' },
    @{ Cell = 'E25'; Value = '  -4.38%  ' },
    @{ Cell = 'D26'; Value = '8.997' },
    @{ Cell = 'D27'; Value = '161.29' },
    @{ Cell = 'E27'; Value = '  -0.87%  ' },
    @{ Cell = 'D28'; Value = '18.74' },
    @{ Cell = 'E28'; Value = '  +1.26%  ' },
    @{ Cell = 'D29'; Value = '1.580' },
    @{ Cell = 'E29'; Value = '  +4.93%  ' },
    @{ Cell = 'D30'; Value = '4.301' },
    @{ Cell = 'E30'; Value = '  -2.27%  ' },
    @{ Cell = 'D31'; Value = '4.258' },
    @{ Cell = 'E31'; Value = '  -0.49%  ' },
    @{ Cell = 'E32'; Value = '  -1.53%  ' },
    @{ Cell = 'D33'; Value = '0.05261' },
    @{ Cell = 'E33'; Value = '  -1.60%  ' },
    @{ Cell = 'D34'; Value = '1.892' },
    @{ Cell = 'E34'; Value = '  -2.18%  ' },
    @{ Cell = 'D35'; Value = '0.7581' },
    @{ Cell = 'E35'; Value = '  -0.25%  ' },
    @{ Cell = 'D36'; Value = '1.180' },
    @{ Cell = 'E36'; Value = '  +0.28%  ' },
    @{ Cell = 'E37'; Value = '  +0.29%  ' },
    @{ Cell = 'B38'; Value = 'Maker' },
    @{ Cell = 'C38'; Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr' },
    @{ Cell = 'D38'; Value = '1.276.36' },
    @{ Cell = 'E38'; Value = '  +1.00%  ' },
    @{ Cell = 'B39'; Value = 'VeChain' },
    @{ Cell = 'C39'; Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet' },
    @{ Cell = 'D39'; Value = '0.01869' },
    @{ Cell = 'E39'; Value = '  +0.01%  ' },
    @{ Cell = 'D40'; Value = '2.757' },
    @{ Cell = 'D41'; Value = '0.8997' },
    @{ Cell = 'E41'; Value = '  -0.48%  ' },
    @{ Cell = 'D42'; Value = '109.77' },
    @{ Cell = 'E42'; Value = '  -2.84%  ' },
    @{ Cell = 'D43'; Value = '5.966' },
    @{ Cell = 'E43'; Value = '  -7.26%  ' },
    @{ Cell = 'D44'; Value = '70.55' },
    @{ Cell = 'E44'; Value = '  -4.69%  ' },
    @{ Cell = 'E45'; Value = '  -0.15%  ' },
    @{ Cell = 'E46'; Value = '  -2.68%  ' },
    @{ Cell = 'D47'; Value = '2.029.49' },
    @{ Cell = 'E47'; Value = '  +0.23%  ' },
    @{ Cell = 'B48'; Value = 'EnergySwap' },
    @{ Cell = 'C48'; Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens' },
    @{ Cell = 'D48'; Value = '9.605' },
    @{ Cell = 'E48'; Value = '  +1.39%  ' },
    @{ Cell = 'B49'; Value = 'RenderToken' },
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr' },
    @{ Cell = 'D49'; Value = '1.798' },
    @{ Cell = 'E49'; Value = '  +0.04%  ' },
    @{ Cell = 'D50'; Value = '0.5174' },
    @{ Cell = 'E50'; Value = '  -0.47%  ' },
    @{ Cell = 'D51'; Value = '0.4298' },
    @{ Cell = 'E51'; Value = '  -0.89%  ' }
)

foreach ($change in $changes) {
    $range = $ws.Range($change.Cell)
    # Preserve the cell's existing style while forcing the assigned text to
    # stay a literal string (these columns hold text such as "243.79" or
    # "7.600" that Excel would otherwise auto-convert to a number and mangle
    # trailing zeros / thousands separators).
    $originalStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $change.Value
    $range.Style = $originalStyle
}
